# feat: add 2022-Q3 data
#
# 1. Duplicate the "2022-Q2" sheet (it already has the exact formatting/
#    text-typed number cells we need as a template) and place the copy
#    right after "2022-Q2"; rename it to "2022-Q3" and move it to sit
#    right after "总计" (i.e. before "2022-Q2").
# 2. Trim the duplicate down to the two funds that make up the Q3 numbers
#    and overwrite the figures.
# 3. Insert the new "2022-Q3" row at the top of the "总计" summary sheet,
#    shifting the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q3" sheet from the "2022-Q2" template
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")

# Copy placed *before* the template == right after "总计".
$template.Copy($template)

$q3temp = $wb.Worksheets.Item("2022-Q2 (2)")
$q3temp.Name = "2022-Q3"

# Re-fetch by name for a stable reference before further edits.
$q3 = $wb.Worksheets.Item("2022-Q3")

# The template had 4 fund rows (009686, 501201, 168401, 009687); Q3 only
# has 2 (009686, 009687) so drop the middle two rows.
$q3.Range("A3:H4").EntireRow.Delete()

# ---------------------------------------------------------------------
# Step 2: overwrite the figures for the two remaining funds
# ---------------------------------------------------------------------
# Row 2: fund 009686
$q3.Range("H2").Value = 2
$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "10.76"
$q3.Range("D2").Style = "Normal"
$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "64.78"
$q3.Range("E2").Style = "Normal"
$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "3.86"
$q3.Range("F2").Style = "Normal"
$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.4153"
$q3.Range("G2").Style = "Normal"

# Row 3: fund 009687
$q3.Range("A3").Value = 1
$q3.Range("H3").Value = 2
$q3.Range("D3").NumberFormat = "@"
$q3.Range("D3").Value = "0.43"
$q3.Range("D3").Style = "Normal"
$q3.Range("E3").NumberFormat = "@"
$q3.Range("E3").Value = "64.78"
$q3.Range("E3").Style = "Normal"
$q3.Range("F3").NumberFormat = "@"
$q3.Range("F3").Value = "3.86"
$q3.Range("F3").Style = "Normal"
$q3.Range("G3").NumberFormat = "@"
$q3.Range("G3").Value = "0.0166"
$q3.Range("G3").Style = "Normal"

# ---------------------------------------------------------------------
# Step 3: add the 2022-Q3 row to the "总计" summary sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Range("A2:D4").Copy($total.Range("A3:D5"))

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.43

# Column A is a running 0-based index - renumber it after the shift.
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
